# Fixing strain number and entering marker info for off by one errors.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix block 1: bioSample rows 10-20 (sheet rows 11-21, first replicate) ---
# Strain numbers were off by one (duplicated TDY2194 / missing TDY2193), and
# marker_1 (NAT) was never entered for these rows. Shift the strain column up
# by one TDY number and fill in the marker_1 column.
$ws.Range("E11").Value = "TDY2193"
$ws.Range("J11").Value = "NAT"

$ws.Range("E12").Value = "TDY2194"
$ws.Range("J12").Value = "NAT"

$ws.Range("E13").Value = "TDY2195"
$ws.Range("J13").Value = "NAT"

$ws.Range("E14").Value = "TDY2196"
$ws.Range("J14").Value = "NAT"

$ws.Range("E15").Value = "TDY2197"
$ws.Range("J15").Value = "NAT"

$ws.Range("E16").Value = "TDY2198"
$ws.Range("J16").Value = "NAT"

$ws.Range("E17").Value = "TDY2199"
$ws.Range("J17").Value = "NAT"

$ws.Range("E18").Value = "TDY2200"
$ws.Range("J18").Value = "NAT"

$ws.Range("E19").Value = "TDY2201"
$ws.Range("J19").Value = "NAT"

$ws.Range("E20").Value = "TDY2202"
$ws.Range("J20").Value = "NAT"

$ws.Range("E21").Value = "TDY2203"
$ws.Range("J21").Value = "NAT"

# --- Fix block 2: bioSample rows 32-42 (sheet rows 33-43, second replicate) ---
# Same off-by-one / missing marker_1 problem repeats for the second batch.
$ws.Range("E33").Value = "TDY2193"
$ws.Range("J33").Value = "NAT"

$ws.Range("E34").Value = "TDY2194"
$ws.Range("J34").Value = "NAT"

$ws.Range("E35").Value = "TDY2195"
$ws.Range("J35").Value = "NAT"

$ws.Range("E36").Value = "TDY2196"
$ws.Range("J36").Value = "NAT"

$ws.Range("E37").Value = "TDY2197"
$ws.Range("J37").Value = "NAT"

$ws.Range("E38").Value = "TDY2198"
$ws.Range("J38").Value = "NAT"

$ws.Range("E39").Value = "TDY2199"
$ws.Range("J39").Value = "NAT"

$ws.Range("E40").Value = "TDY2200"
$ws.Range("J40").Value = "NAT"

$ws.Range("E41").Value = "TDY2201"
$ws.Range("J41").Value = "NAT"

$ws.Range("E42").Value = "TDY2202"
$ws.Range("J42").Value = "NAT"

$ws.Range("E43").Value = "TDY2203"
$ws.Range("J43").Value = "NAT"

# --- Drop the now-unused marker_2 column (column K) entirely ---
$ws.Columns("K").Delete()

# Reflect where the editor's cursor ended up after making these edits.
$ws.Range("D15").Select()
